# Adds the "rowslice" and "colslice" filter examples to the df_filters sheet,
# mirroring the upstream commit "added colslice and rowslice filters (#1645)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("df_filters")

# --- New filter expression labels in column A (write these first so the
#     shared-string table assigns them indices 44 and 45, ahead of the two
#     whitespace-only strings below, matching the target workbook). ---
$ws.Range("A116").Value = "{{ df2 | rowslice(0, 3) }}"
$ws.Range("A122").Value = "{{ df2 | colslice(3) | rowslice(0, 2) }}"

# --- Existing example rows 20 and 25 gain a couple of whitespace-only
#     shared strings (single-space / double-space placeholders). ---
$ws.Range("I20").Value = " "
$ws.Range("J20").Value = "  "
$ws.Range("H25").Value = " "

# --- Give the new G115:K124 block the same shaded-cell formatting used by
#     the rest of the sheet's data blocks (style index "2" / light fill),
#     by copying the format from an existing ten-row block further up. ---
$ws.Range("G104:K113").Copy()
$ws.Range("G115:K124").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- rowslice(0, 3) example block (rows 116-119: header + 3 data rows) ---
$ws.Range("H116").Value = "name"
$ws.Range("I116").Value = "b"
$ws.Range("J116").Value = "c"
$ws.Range("K116").Value = "d"

$ws.Range("G117").Value = 0
$ws.Range("H117").Value = "a"
$ws.Range("I117").Value = 4
$ws.Range("J117").Value = 1
$ws.Range("K117").Value = 1

$ws.Range("G118").Value = 1
$ws.Range("H118").Value = "b"
$ws.Range("I118").Value = 2
$ws.Range("J118").Value = 2
$ws.Range("K118").Value = 1

$ws.Range("G119").Value = 2
$ws.Range("H119").Value = "c"
$ws.Range("I119").Value = 6
$ws.Range("J119").Value = 5
$ws.Range("K119").Value = 1

# Rows 115, 120 and 121 stay blank (only the shaded formatting applies).

# --- colslice(3) | rowslice(0, 2) example block (rows 122-124) ---
$ws.Range("H122").Value = "d"

$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 1

$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 1

# --- Update the sheet's view: drop the old scrolled/selected cell and
#     select G16 instead (matches the new sheetView in the diff). ---
$ws.Activate()
$ws.Range("G16").Select()

# --- Widen the saved workbook window (best effort; cosmetic UI metadata). ---
try {
  $excel.ActiveWindow.Width = 38400
} catch {
}
